$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "APPLE iPhone SE (White, 128 GB)"

# Rewrite every row's value into column A (instead of the old diagonal
# layout A1, B2, C3, ... AE31) and update the product description text.
for ($row = 1; $row -le 31; $row++) {
    $ws.Cells.Item($row, 1).Value2 = $newValue
}

# Clear out the old diagonal cells (B2, C3, ..., AE31) now that their
# values live in column A.
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, $row).ClearContents()
}
